$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rules")
$ws2 = $wb.Worksheets.Item("Hoja1")

# --- Sheet "Rules": clear the old suggestCards header block (rows 4:18) ---
$ws1.Range("B4:D4").UnMerge()
$ws1.Range("B4:D18").ClearContents()

# --- Sheet "Rules": clear the old suggestCards3 summary block (rows 141:167) ---
$ws1.Range("B141:E141").UnMerge()
$ws1.Range("B149:B151").UnMerge()
$ws1.Range("C149:C151").UnMerge()
$ws1.Range("D149:D151").UnMerge()
$ws1.Range("B155:B158").UnMerge()
$ws1.Range("C155:C158").UnMerge()
$ws1.Range("D155:D158").UnMerge()
$ws1.Range("B159:B160").UnMerge()
$ws1.Range("C159:C160").UnMerge()
$ws1.Range("D159:D160").UnMerge()
$ws1.Range("B161:B162").UnMerge()
$ws1.Range("C161:C162").UnMerge()
$ws1.Range("D161:D162").UnMerge()
$ws1.Range("B163:B165").UnMerge()
$ws1.Range("C163:C165").UnMerge()
$ws1.Range("D163:D165").UnMerge()
$ws1.Range("B141:E167").ClearContents()

# --- Sheet "Hoja1": rename the rule title from suggestCards3(..) to suggestCards(..) ---
$ws2.Range("B2").Value = "Rules String[] suggestCards (String situation, Double income)"
